$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "QZMWMK"
$ws.Range("B66").Value = "Fil de fusor del medio HP/CANON"
$ws.Range("C66").Value = "HP M252 M277 M377 M477 M452 M454 M455 M479, Canon MF 732 734 735"
$ws.Range("D66").Value = 40000
$ws.Range("E66").Value = 150000
$ws.Range("F66").Value = 10
$ws.Range("G66").Value = 3
$ws.Range("H66").Formula = "=(E66-D66)*G66"
$ws.Range("I66").Formula = "=D66*F66"
$ws.Range("J66").Value = 400000
